$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-54 down to 5-55
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44921
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 100112028
$ws.Range("G4").Value = "Sandia"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 840
$ws.Range("K4").Value = 430
$ws.Range("L4").Value = 450
$ws.Range("M4").Value = 438
$ws.Range("N4").Value = "$/kilo (volumen en unidades)"
$ws.Range("O4").Value = "Perú"
$ws.Range("P4").Value = 438
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Hortaliza"
